$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - new numeric values in E1, F1
$ws.Range("E1").Value = 0.08
$ws.Range("F1").Value = 0.06

# Row 2 - header labels (reuse existing shared string "factor_correction")
$ws.Range("E2").Value = "factor_correction"
$ws.Range("F2").Value = "factor_correction"

# Rows 3-15 - ratio formulas in columns E and F
$ws.Range("E3").Formula = "=56.5051/67.3424"
$ws.Range("F3").Formula = "=56.5051/67.3424"

$ws.Range("E4").Formula = "=51.6663/64.7336"
$ws.Range("F4").Formula = "=51.6663/64.7336"

$ws.Range("E5").Formula = "=42.9408/47.1749"
$ws.Range("F5").Formula = "=42.9408/47.1749"

$ws.Range("E6").Formula = "=31.7833/35.5311"
$ws.Range("F6").Formula = "=31.7833/35.5311"

$ws.Range("E7").Formula = "=51.591/56.6322"
$ws.Range("F7").Formula = "=51.5901/56.6322"

$ws.Range("E8").Formula = "=38.552/45.678"
$ws.Range("F8").Formula = "=38.552/45.678"

$ws.Range("E9").Formula = "=38.6/41.1682"
$ws.Range("F9").Formula = "=38.6/41.16822"

$ws.Range("E10").Formula = "=34.4685/37.0638"
$ws.Range("F10").Formula = "=34.468572/37.063807"

$ws.Range("E11").Formula = "=28.6625/29.4184"
$ws.Range("F11").Formula = "=28.6625/29.4184"

$ws.Range("E12").Formula = "=16.5524/16.7974"
$ws.Range("F12").Formula = "=16.5524/16.7974"

$ws.Range("E13").Formula = "=29.2376/29.9596"
$ws.Range("F13").Formula = "=29.2376/29.9596"

$ws.Range("E14").Formula = "=30.5607/31.0304"
$ws.Range("F14").Formula = "=30.5607/31.0304"

$ws.Range("E15").Formula = "=31.4635/34.6892"
$ws.Range("F15").Formula = "=31.46358/34.6892"

# Column widths for E:F to match bestFit width used by column C (16.28515625)
$ws.Columns("E:F").EntireColumn.AutoFit() | Out-Null

# Update selection to match target (E16)
$ws.Range("E16").Select()

$wb.Save()
